# Added new in-game tunes: "Killed by Shadow" and "Meet Jaffar"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Style bookkeeping -----------------------------------------------------
# The old "Good" (green, centered) cell style is being retired; its cellXfs
# slot (index 6) gets reused for a plain style with left alignment instead.
# Remove the built-in "Good" style, which also resets any cell still using
# it (G5) back to the default "Normal" style.
$wb.Styles("Good").Delete()

# --- Row 5 (s_Danger): G5 loses the old "Good" look, becomes centered ------
$ws.Range("G5").Style = "Normal"
$ws.Range("G5").HorizontalAlignment = -4108   # xlCenter

# --- Row 6 (s_Sword) --------------------------------------------------------
# H6 loses its yellow highlight (back to default/Normal formatting).
$ws.Range("H6").Style = "Normal"
# J6 filename is replaced by the new consolidated tune name and picks up the
# new plain left-aligned look (reusing cellXfs slot 6).
$ws.Range("J6").Value = "m-killguard-or-sword"
$ws.Range("J6").Style = "Normal"
$ws.Range("J6").HorizontalAlignment = -4131   # xlLeft

# --- Row 8 (s_Shadow): "Killed by Shadow" tune -----------------------------
$ws.Range("H8").Style = "Normal"
$ws.Range("H8").Value = "pop_music_shadow"
$ws.Range("J8").Value = "m-e4-killedbyshadow"

# --- Row 9 (s_Vict): same consolidated filename as row 6 -------------------
$ws.Range("J9").Value = "m-killguard-or-sword"
$ws.Range("J9").Style = "Normal"
$ws.Range("J9").HorizontalAlignment = -4131   # xlLeft

# --- Row 10 (s_Stairs): code label cleared, highlight kept -----------------
$ws.Range("H10").ClearContents()

# --- Row 12 (s_Jaffar): "Meet Jaffar" tune ---------------------------------
$ws.Range("H12").Value = "pop_music_jaffar"
$ws.Range("J12").Value = "m_jaffar"
$ws.Range("J12").Style = "Normal"
$ws.Range("J12").HorizontalAlignment = -4131  # xlLeft

# --- Scroll the window so column G is the first visible column -------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
